$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$errMsg = "Failed to select the date: Message: no such element: Unable to locate element: {`"method`":`"css selector`",`"selector`":`"#restProfileSideBarDtpDayPicker-label-wrapper button[aria-label*='10-01']`"}`n  (Session info: chrome=128.0.6613.120); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception`nStacktrace:`n`tGetHandleVerifier [0x00007FF7F082B5D2+29090]`n`t(No symbol) [0x00007FF7F079E689]`n`t(No symbol) [0x00007FF7F065B1CA]`n`t(No symbol) [0x00007FF7F06AEFD7]`n`t(No symbol) [0x00007FF7F06AF22C]`n`t(No symbol) [0x00007FF7F06F97F7]`n`t(No symbol) [0x00007FF7F06D672F]`n`t(No symbol) [0x00007FF7F06F65D9]`n`t(No symbol) [0x00007FF7F06D6493]`n`t(No symbol) [0x00007FF7F06A09B1]`n`t(No symbol) [0x00007FF7F06A1B11]`n`tGetHandleVerifier [0x00007FF7F0B48C5D+3295277]`n`tGetHandleVerifier [0x00007FF7F0B94843+3605523]`n`tGetHandleVerifier [0x00007FF7F0B8A707+3564247]`n`tGetHandleVerifier [0x00007FF7F08E6EB6+797318]`n`t(No symbol) [0x00007FF7F07A980F]`n`t(No symbol) [0x00007FF7F07A53F4]`n`t(No symbol) [0x00007FF7F07A5580]`n`t(No symbol) [0x00007FF7F0794A1F]`n`tBaseThreadInitThunk [0x00007FFC979C257D+29]`n`tRtlUserThreadStart [0x00007FFC9896AF28+40]`n"

$rows = @(
    @{ r = 10; a = "2024-09-09 17:12:14"; f = "17:12:14" },
    @{ r = 11; a = "2024-09-09 17:12:35"; f = "17:12:35" },
    @{ r = 12; a = "2024-09-09 17:12:56"; f = "17:12:56" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 2).Value = "monitor_availability"
    $ws.Cells.Item($r, 3).Value = "https://www.opentable.com/r/bar-spero-washington/"
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = "2024-09-09"
    $ws.Cells.Item($r, 5).ClearFormats()
    $ws.Cells.Item($r, 6).Value = $row.f
    $ws.Cells.Item($r, 4).Value = $errMsg
}
